$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first worksheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 540
$wsExhibition.Range("F7").Value = 2633
$wsExhibition.Range("F9").Value = 7234
$wsExhibition.Range("F13").Value = 170

# Sheet "全部类型" (All Types) - fourth worksheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 540
$wsAll.Range("F9").Value = 2633
$wsAll.Range("F11").Value = 7234
$wsAll.Range("F17").Value = 170
